$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Satisfaction between Age Groups")

# Rename the second "Tests" header row (A13:F13) from TestA..TestE to Test1..Test5,
# matching the naming used by the other tables on this sheet.
$ws.Range("B13").Value = "Test1"
$ws.Range("C13").Value = "Test2"
$ws.Range("D13").Value = "Test3"
$ws.Range("E13").Value = "Test4"
$ws.Range("F13").Value = "Test5"

# Updated per-user results (rows 14:16) reflecting the new test emissions.
$ws.Range("B14").Value = 4
$ws.Range("C14").Value = 6
$ws.Range("D14").Value = 4
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 10

$ws.Range("B15").Value = 10
$ws.Range("C15").Value = 4
$ws.Range("D15").Value = 5
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 10

$ws.Range("B16").Value = 6
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 10

# Rows 19:21 (Min/Average/Max) recalc automatically from the formulas already in place.

# Move the active selection as recorded in the saved view state.
$ws.Range("A24").Select()
